# Applies the "balancos concatenados" correction pass: a handful of
# floating point roundings get nudged to their newly-recomputed values,
# and the empty "Part. de Acionistas Nao Controladores" row (78) has its
# trailing zero-filled columns (P:AH) turned into real blank cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58 - Receitas da Intermediacao Financeira
$ws.Range("H58").Value  = 229406.944
$ws.Range("L58").Value  = 358572.992
$ws.Range("P58").Value  = 456922.016
$ws.Range("T58").Value  = 46166.08
$ws.Range("X58").Value  = 615957.12
$ws.Range("AF58").Value = 510714.048

# Row 59 - Despesas da Intermediacao Financeira
$ws.Range("H59").Value  = -114373.008
$ws.Range("L59").Value  = -156210.96
$ws.Range("T59").Value  = 108250.016
$ws.Range("X59").Value  = -418814.912
$ws.Range("AB59").Value = -472089.056

# Row 60 - Resultado Bruto Intermediacao Financeira
$ws.Range("H60").Value  = 115034.016
$ws.Range("L60").Value  = 202362
$ws.Range("T60").Value  = 154415.984
$ws.Range("AF60").Value = -213181.984

# Row 61 - Outras Despesas/Receitas Operacionais
$ws.Range("L61").Value  = -86589
$ws.Range("P61").Value  = -58017
$ws.Range("T61").Value  = -96830.984
$ws.Range("X61").Value  = -185326.992
$ws.Range("AB61").Value = -107957.992

# Row 63 - Despesas de Pessoal
$ws.Range("AF63").Value = -53299

# Row 67 - Outras Receitas Operacionais
$ws.Range("T67").Value  = -28577.992
$ws.Range("AF67").Value = -68225

# Row 69 - Resultado Operacional
$ws.Range("L69").Value  = 115773
$ws.Range("P69").Value  = 133854.008
$ws.Range("T69").Value  = 57585.016
$ws.Range("AF69").Value = -363140.96

# Row 73 - Receitas
$ws.Range("H73").Value = 4238
$ws.Range("P73").Value = 133984.024

# Row 74 - Despesas
$ws.Range("H74").Value = -50303
$ws.Range("L74").Value = 17424.008
$ws.Range("P74").Value = 17603

# Row 78 - Part. de Acionistas Nao Controladores: columns P through AH
# were all zero placeholders; they become genuinely blank cells now.
$ws.Range("P78:AH78").ClearContents()

# Row 79 - Provisao para IR e Contribuicao Social
$ws.Range("H79").Value = 20282.008
$ws.Range("L79").Value = 85845
$ws.Range("P79").Value = 69474.008
$ws.Range("T79").Value = 30634.008
